$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$updates = @{
  "H40" = 7272
  "I40" = 5999.2
  "J40" = 8332.666999999999
  "K40" = 5999.2
  "L40" = 8332.666999999999
  "M40" = -5824.2
  "N40" = -8682.666999999999
  "H55" = 1181.4
  "I55" = 309.2
  "J55" = 2053.6
  "K55" = 309.2
  "L55" = 2053.6
  "M55" = -95.19999999999999
  "N55" = -2481.6
  "H69" = 7384.316
  "J69" = 7776.5884
  "L69" = 23329.7652
  "N69" = -25077.7652
  "H72" = 7384.316
  "J72" = 7776.5884
  "L72" = 69989.2956
  "N72" = -78725.2956
  "H80" = 327.92307
  "I80" = 207.33333
  "K80" = 621.99999
  "M80" = 376.00001
  "H83" = 327.92307
  "I83" = 207.33333
  "K83" = 1865.99997
  "M83" = 3126.00003
  "H87" = 95354
  "J87" = 95354
  "L87" = 95354
  "N87" = -97850
  "H90" = 95354
  "J90" = 95354
  "L90" = 286062
  "N90" = -298542
  "H92" = 224.83333
  "J92" = 216.33333
  "L92" = 216.33333
  "N92" = -2712.33333
  "H130" = 50000
  "J130" = 50000
  "L130" = 50000
  "N130" = -60040
}
foreach ($key in $updates.Keys) {
  $ws.Range($key).Value = $updates[$key]
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$updates = @{
  "H2" = 1249
  "I2" = 1293.1111
  "K2" = 1293.1111
  "M2" = -1180.1111
  "H32" = 20002266
  "I32" = 2831.75
  "K32" = 2831.75
  "M32" = -2544.75
  "H45" = 2942.5833
  "I45" = 2245.875
  "K45" = 2245.875
  "M45" = -1868.875
  "H97" = 1412.625
  "I97" = 900.2857
  "K97" = 900.2857
  "M97" = -404.2857
  "H102" = 7816901.5
  "I102" = 13890234
  "J102" = 8331.286
  "K102" = 13890234
  "L102" = 8331.286
  "M102" = -13888612
  "N102" = -11575.286
  "H116" = 1249
  "I116" = 1293.1111
  "K116" = 1293.1111
  "M116" = 1000.8889
  "H132" = 2279.9333
  "I132" = 2018.25
  "J132" = 3326.6667
  "K132" = 6054.75
  "L132" = 9980.000100000001
  "M132" = -3524.75
  "N132" = -15040.0001
}
foreach ($key in $updates.Keys) {
  $ws.Range($key).Value = $updates[$key]
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$updates = @{
  "H3" = 1249
  "I3" = 1293.1111
  "K3" = 1293.1111
  "M3" = -1179.1111
  "H20" = 3056.75
  "I20" = 1614
  "K20" = 1614
  "M20" = -1367
  "H81" = 45476.11
  "J81" = 45476.11
  "L81" = 45476.11
  "N81" = -47598.11
  "H84" = 45476.11
  "J84" = 45476.11
  "L84" = 136428.33
  "N84" = -147036.33
  "H86" = 3730.3684
  "I86" = 1636.7693
  "K86" = 1636.7693
  "M86" = -513.7692999999999
  "H89" = 3730.3684
  "I89" = 1636.7693
  "K89" = 8183.8465
  "M89" = -2567.8465
  "H94" = 652.1818
  "I94" = 692.4
  "J94" = 250
  "K94" = 692.4
  "L94" = 250
  "M94" = -241.4
  "N94" = -1152
  "H99" = 250001500
  "I99" = 333334660
  "K99" = 333334660
  "M99" = -333333162
  "H105" = 4133721
  "I105" = 5683099
  "K105" = 5683099
  "M105" = -5681352
  "H107" = 45460536
  "I107" = 100001180
  "K107" = 100001180
  "M107" = -99999260
}
foreach ($key in $updates.Keys) {
  $ws.Range($key).Value = $updates[$key]
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$updates = @{
  "H105" = 2477.3333
  "I105" = 880.8333
  "K105" = 880.8333
  "M105" = 866.1667
  "H107" = 1227.1111
  "I107" = 466
  "J107" = 2423.1428
  "K107" = 466
  "L107" = 2423.1428
  "M107" = 1454
  "N107" = -6263.1428
  "H132" = 2490.7856
  "I132" = 2490.7856
  "K132" = 7472.3568
  "M132" = -4942.3568
  "H134" = 1414.3334
  "J134" = 597.5
  "L134" = 1792.5
  "N134" = -6862.5
}
foreach ($key in $updates.Keys) {
  $ws.Range($key).Value = $updates[$key]
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$updates = @{
  "H15" = 147.71428
  "I15" = 97
  "J15" = 168
  "K15" = 291
  "L15" = 504
  "M15" = -151
  "N15" = -784
  "H22" = 175
  "J22" = 300
  "L22" = 900
  "N22" = -1238
  "H27" = 175
  "J27" = 300
  "L27" = 900
  "N27" = -1104
  "H34" = 2313.5454
  "I34" = 175.75
  "J34" = 3535.1428
  "K34" = 527.25
  "L34" = 10605.4284
  "M34" = -443.25
  "N34" = -10773.4284
  "H75" = 5039.6
  "J75" = 6216.25
  "L75" = 18648.75
  "N75" = -20644.75
  "H78" = 5039.6
  "J78" = 6216.25
  "L78" = 55946.25
  "N78" = -65930.25
  "H92" = 1811.1111
  "I92" = 1533.3334
  "J92" = 2366.6667
  "K92" = 4600.0002
  "L92" = 7100.000100000001
  "M92" = -3352.0002
  "N92" = -9596.000100000001
  "H124" = 1097.3334
  "J124" = 1109.4286
  "L124" = 3328.2858
  "N124" = -13148.2858
}
foreach ($key in $updates.Keys) {
  $ws.Range($key).Value = $updates[$key]
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$updates = @{
  "H63" = 24999.334
  "J63" = 24999.334
  "L63" = 24999.334
  "N63" = -26371.334
  "H66" = 24999.334
  "J66" = 24999.334
  "L66" = 74998.00199999999
  "N66" = -81862.00199999999
  "H70" = 10000
  "I70" = 0
  "K70" = 0
  "H73" = 10000
  "I73" = 0
  "K73" = 0
  "H107" = 491.14285
  "I107" = 287.8
  "K107" = 287.8
  "M107" = 1632.2
  "H132" = 0
  "I132" = 0
  "K132" = 0
}
foreach ($key in $updates.Keys) {
  $ws.Range($key).Value = $updates[$key]
}
$ws.Range("M70").ClearContents()
$ws.Range("M73").ClearContents()
$ws.Range("M132").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$updates = @{
  "H22" = 1039.6875
  "I22" = 775
  "J22" = 1198.5
  "K22" = 775
  "L22" = 1198.5
  "M22" = -480
  "N22" = -1788.5
  "H27" = 1039.6875
  "I27" = 775
  "J27" = 1198.5
  "K27" = 775
  "L27" = 1198.5
  "M27" = -668
  "N27" = -1412.5
  "H69" = 51990
  "J69" = 51990
  "L69" = 51990
  "N69" = -53612
  "H72" = 51990
  "J72" = 51990
  "L72" = 155970
  "N72" = -164082
  "H100" = 7273
  "I100" = 4400.6
  "J100" = 9666.666999999999
  "K100" = 4400.6
  "L100" = 9666.666999999999
  "M100" = -3859.6
}
foreach ($key in $updates.Keys) {
  $ws.Range($key).Value = $updates[$key]
}
